$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object 'object[,]' 1,16
$row2[0,0] = 3
$row2[0,1] = 1
$row2[0,2] = 103.4275383333333
$row2[0,3] = 310.282615
$row2[0,4] = 0.2485530285127421
$row2[0,5] = 0.2485530285127421
$row2[0,6] = 3
$row2[0,7] = 1
$row2[0,8] = 19.655164
$row2[0,9] = 58.965492
$row2[0,10] = 0.2909311988313945
$row2[0,11] = 0.2909311988313944
$row2[0,12] = 2032.885228057954
$row2[0,13] = 18295.96705252158
$row2[0,14] = 0.07231183055838583
$row2[0,15] = 0.0723118305583858
$ws.Range("E2:T2").Value2 = $row2

$row3 = New-Object 'object[,]' 1,16
$row3[0,0] = 3
$row3[0,1] = 1
$row3[0,2] = 103.4275383333333
$row3[0,3] = 310.282615
$row3[0,4] = 0.2485530285127421
$row3[0,5] = 0.2485530285127421
$row3[0,6] = 3
$row3[0,7] = 1
$row3[0,8] = 1.429363
$row3[0,9] = 4.288089
$row3[0,10] = 0.02115710106286768
$row3[0,11] = 0.02115710106286767
$row3[0,12] = 147.8354964747484
$row3[0,13] = 1330.519468272735
$row3[0,14] = 0.005258661543725916
$row3[0,15] = 0.005258661543725915
$ws.Range("E3:T3").Value2 = $row3

$row4 = New-Object 'object[,]' 1,16
$row4[0,0] = 3
$row4[0,1] = 1
$row4[0,2] = 103.4275383333333
$row4[0,3] = 310.282615
$row4[0,4] = 0.2485530285127421
$row4[0,5] = 0.2485530285127421
$row4[0,6] = 3
$row4[0,7] = 1
$row4[0,8] = 22.55245966666666
$row4[0,9] = 67.65737899999999
$row4[0,10] = 0.3338163002567673
$row4[0,11] = 0.3338163002567672
$row4[0,12] = 2332.54538668512
$row4[0,13] = 20992.90848016608
$row4[0,14] = 0.08297105239573836
$row4[0,15] = 0.08297105239573833
$ws.Range("E4:T4").Value2 = $row4

$row5 = New-Object 'object[,]' 1,16
$row5[0,0] = 3
$row5[0,1] = 1
$row5[0,2] = 103.4275383333333
$row5[0,3] = 310.282615
$row5[0,4] = 0.2485530285127421
$row5[0,5] = 0.2485530285127421
$row5[0,6] = 3
$row5[0,7] = 1
$row5[0,8] = 23.92250533333333
$row5[0,9] = 71.767516
$row5[0,10] = 0.3540953998489707
$row5[0,11] = 0.3540953998489706
$row5[0,12] = 2474.245837392705
$row5[0,13] = 22268.21253653434
$row5[0,14] = 0.08801148401489203
$row5[0,15] = 0.08801148401489201
$ws.Range("E5:T5").Value2 = $row5

$row6 = New-Object 'object[,]' 1,16
$row6[0,0] = 3
$row6[0,1] = 1
$row6[0,2] = 216.130539
$row6[0,3] = 648.391617
$row6[0,4] = 0.5193964865470273
$row6[0,5] = 0.5193964865470272
$row6[0,6] = 3
$row6[0,7] = 1
$row6[0,8] = 19.655164
$row6[0,9] = 58.965492
$row6[0,10] = 0.2909311988313945
$row6[0,11] = 0.2909311988313944
$row6[0,12] = 4248.081189453395
$row6[0,13] = 38232.73070508056
$row6[0,14] = 0.1511086424999409
$row6[0,15] = 0.1511086424999408
$ws.Range("E6:T6").Value2 = $row6

$row7 = New-Object 'object[,]' 1,16
$row7[0,0] = 3
$row7[0,1] = 1
$row7[0,2] = 216.130539
$row7[0,3] = 648.391617
$row7[0,4] = 0.5193964865470273
$row7[0,5] = 0.5193964865470272
$row7[0,6] = 3
$row7[0,7] = 1
$row7[0,8] = 1.429363
$row7[0,9] = 4.288089
$row7[0,10] = 0.02115710106286768
$row7[0,11] = 0.02115710106286767
$row7[0,12] = 308.928995616657
$row7[0,13] = 2780.360960549913
$row7[0,14] = 0.01098892395757385
$row7[0,15] = 0.01098892395757384
$ws.Range("E7:T7").Value2 = $row7

$row8 = New-Object 'object[,]' 1,16
$row8[0,0] = 3
$row8[0,1] = 1
$row8[0,2] = 216.130539
$row8[0,3] = 648.391617
$row8[0,4] = 0.5193964865470273
$row8[0,5] = 0.5193964865470272
$row8[0,6] = 3
$row8[0,7] = 1
$row8[0,8] = 22.55245966666666
$row8[0,9] = 67.65737899999999
$row8[0,10] = 0.3338163002567673
$row8[0,11] = 0.3338163002567672
$row8[0,12] = 4874.275263532426
$row8[0,13] = 43868.47737179184
$row8[0,14] = 0.1733830135054925
$row8[0,15] = 0.1733830135054924
$ws.Range("E8:T8").Value2 = $row8

$row9 = New-Object 'object[,]' 1,16
$row9[0,0] = 3
$row9[0,1] = 1
$row9[0,2] = 216.130539
$row9[0,3] = 648.391617
$row9[0,4] = 0.5193964865470273
$row9[0,5] = 0.5193964865470272
$row9[0,6] = 3
$row9[0,7] = 1
$row9[0,8] = 23.92250533333333
$row9[0,9] = 71.767516
$row9[0,10] = 0.3540953998489707
$row9[0,11] = 0.3540953998489706
$row9[0,12] = 5170.383971923708
$row9[0,13] = 46533.45574731338
$row9[0,14] = 0.1839159065840202
$row9[0,15] = 0.1839159065840201
$ws.Range("E9:T9").Value2 = $row9

$row10 = New-Object 'object[,]' 1,16
$row10[0,0] = 3
$row10[0,1] = 1
$row10[0,2] = 71.607325
$row10[0,3] = 214.821975
$row10[0,4] = 0.1720839321833696
$row10[0,5] = 0.1720839321833696
$row10[0,6] = 3
$row10[0,7] = 1
$row10[0,8] = 19.655164
$row10[0,9] = 58.965492
$row10[0,10] = 0.2909311988313945
$row10[0,11] = 0.2909311988313944
$row10[0,12] = 1407.4537164763
$row10[0,13] = 12667.0834482867
$row10[0,14] = 0.05006458468972809
$row10[0,15] = 0.05006458468972808
$ws.Range("E10:T10").Value2 = $row10

$row11 = New-Object 'object[,]' 1,16
$row11[0,0] = 3
$row11[0,1] = 1
$row11[0,2] = 71.607325
$row11[0,3] = 214.821975
$row11[0,4] = 0.1720839321833696
$row11[0,5] = 0.1720839321833696
$row11[0,6] = 3
$row11[0,7] = 1
$row11[0,8] = 1.429363
$row11[0,9] = 4.288089
$row11[0,10] = 0.02115710106286768
$row11[0,11] = 0.02115710106286767
$row11[0,12] = 102.352860883975
$row11[0,13] = 921.1757479557751
$row11[0,14] = 0.003640797144499218
$row11[0,15] = 0.003640797144499217
$ws.Range("E11:T11").Value2 = $row11

$row12 = New-Object 'object[,]' 1,16
$row12[0,0] = 3
$row12[0,1] = 1
$row12[0,2] = 71.607325
$row12[0,3] = 214.821975
$row12[0,4] = 0.1720839321833696
$row12[0,5] = 0.1720839321833696
$row12[0,6] = 3
$row12[0,7] = 1
$row12[0,8] = 22.55245966666666
$row12[0,9] = 67.65737899999999
$row12[0,10] = 0.3338163002567673
$row12[0,11] = 0.3338163002567672
$row12[0,12] = 1614.921308900392
$row12[0,13] = 14534.29178010352
$row12[0,14] = 0.05744442157508887
$row12[0,15] = 0.05744442157508886
$ws.Range("E12:T12").Value2 = $row12

$row13 = New-Object 'object[,]' 1,16
$row13[0,0] = 3
$row13[0,1] = 1
$row13[0,2] = 71.607325
$row13[0,3] = 214.821975
$row13[0,4] = 0.1720839321833696
$row13[0,5] = 0.1720839321833696
$row13[0,6] = 3
$row13[0,7] = 1
$row13[0,8] = 23.92250533333333
$row13[0,9] = 71.767516
$row13[0,10] = 0.3540953998489707
$row13[0,11] = 0.3540953998489706
$row13[0,12] = 1713.026614218233
$row13[0,13] = 15417.2395279641
$row13[0,14] = 0.0609341287740534
$row13[0,15] = 0.06093412877405339
$ws.Range("E13:T13").Value2 = $row13

$row14 = New-Object 'object[,]' 1,16
$row14[0,0] = 3
$row14[0,1] = 1
$row14[0,2] = 24.953198
$row14[0,3] = 74.859594
$row14[0,4] = 0.05996655275686102
$row14[0,5] = 0.05996655275686102
$row14[0,6] = 3
$row14[0,7] = 1
$row14[0,8] = 19.655164
$row14[0,9] = 58.965492
$row14[0,10] = 0.2909311988313945
$row14[0,11] = 0.2909311988313944
$row14[0,12] = 490.459199014472
$row14[0,13] = 4414.132791130248
$row14[0,14] = 0.01744614108333964
$row14[0,15] = 0.01744614108333964
$ws.Range("E14:T14").Value2 = $row14

$row15 = New-Object 'object[,]' 1,16
$row15[0,0] = 3
$row15[0,1] = 1
$row15[0,2] = 24.953198
$row15[0,3] = 74.859594
$row15[0,4] = 0.05996655275686102
$row15[0,5] = 0.05996655275686102
$row15[0,6] = 3
$row15[0,7] = 1
$row15[0,8] = 1.429363
$row15[0,9] = 4.288089
$row15[0,10] = 0.02115710106286768
$row15[0,11] = 0.02115710106286767
$row15[0,12] = 35.66717795287401
$row15[0,13] = 321.004601575866
$row15[0,14] = 0.001268718417068695
$row15[0,15] = 0.001268718417068695
$ws.Range("E15:T15").Value2 = $row15

$row16 = New-Object 'object[,]' 1,16
$row16[0,0] = 3
$row16[0,1] = 1
$row16[0,2] = 24.953198
$row16[0,3] = 74.859594
$row16[0,4] = 0.05996655275686102
$row16[0,5] = 0.05996655275686102
$row16[0,6] = 3
$row16[0,7] = 1
$row16[0,8] = 22.55245966666666
$row16[0,9] = 67.65737899999999
$row16[0,10] = 0.3338163002567673
$row16[0,11] = 0.3338163002567672
$row16[0,12] = 562.7559914493472
$row16[0,13] = 5064.803923044125
$row16[0,14] = 0.02001781278044759
$row16[0,15] = 0.02001781278044759
$ws.Range("E16:T16").Value2 = $row16

$row17 = New-Object 'object[,]' 1,16
$row17[0,0] = 3
$row17[0,1] = 1
$row17[0,2] = 24.953198
$row17[0,3] = 74.859594
$row17[0,4] = 0.05996655275686102
$row17[0,5] = 0.05996655275686102
$row17[0,6] = 3
$row17[0,7] = 1
$row17[0,8] = 23.92250533333333
$row17[0,9] = 71.767516
$row17[0,10] = 0.3540953998489707
$row17[0,11] = 0.3540953998489706
$row17[0,12] = 596.9430122387226
$row17[0,13] = 5372.487110148504
$row17[0,14] = 0.0212338804760051
$row17[0,15] = 0.0212338804760051
$ws.Range("E17:T17").Value2 = $row17